# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 11
$ws1.Range("F3").Value = 167
$ws1.Range("F4").Value = 747
$ws1.Range("F5").Value = 64

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 11
$ws4.Range("F4").Value = 167
$ws4.Range("F5").Value = 747
$ws4.Range("F6").Value = 64
